$d = $word.ActiveDocument

# --- Locate the end of the "Database" bullet's text (just before the
# hidden _GoBack bookmark that Word keeps at the very end of the story) ---
$anchor = $d.Content
$found = $anchor.Find.Execute("keeps other people out.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find end of Database paragraph"
}
$anchor.Collapse(0)

# Insert the new "Food" bullet as its own paragraph (inherits the
# ListParagraph / ilvl=1 / numId=1 formatting from the paragraph it is
# split out of, exactly like pressing Enter at that point in Word),
# followed by a temporary extra paragraph ("~MARKER~") that we use to
# safely re-seat the _GoBack bookmark afterwards.
$anchor.InsertAfter("`rFood " + [char]0x2013 + " The application allows for party-goers to order food. Once again these open tabs to popular places so the user is in control of their private information. `r~MARKER~")

# --- Re-seat the _GoBack bookmark at the end of the new "Food" bullet ---
# Grab a safe (non-degenerate-at-end) position: the start of the
# temporary marker paragraph, which sits right after the "Food" bullet.
$lastPara = $d.Paragraphs.Last
$safePos = $lastPara.Range.Start
$bmRange = $d.Range($safePos, $safePos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the marker paragraph's text, then merge its now-empty paragraph
# back into the "Food" bullet by deleting the paragraph mark that
# precedes the bookmark. This leaves the bookmark collapsed immediately
# after "information. " with no leftover paragraph.
$markerText = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$markerText.Delete()
$mergeMark = $d.Range($safePos - 1, $safePos)
$mergeMark.Delete()
